$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 108 (hunk 0)
$ws.Range("H108").Value = 68000
$ws.Range("J108").Value = 68000
$ws.Range("L108").Value = 68000
$ws.Range("N108").Value = -75680

# Row 116 (hunk 1)
$ws.Range("H116").Value = 36033.332
$ws.Range("I116").Value = 100000
$ws.Range("J116").Value = 4050
$ws.Range("K116").Value = 100000
$ws.Range("L116").Value = 4050
$ws.Range("M116").Value = -96558
$ws.Range("N116").Value = -10934

# Row 132 (hunk 2)
$ws.Range("H132").Value = 866.9722
$ws.Range("I132").Value = 857.4857
$ws.Range("J132").Value = 1199
$ws.Range("K132").Value = 2572.4571
$ws.Range("L132").Value = 3597
$ws.Range("M132").Value = -42.45709999999963
$ws.Range("N132").Value = -8657

# Row 139 (hunk 3)
$ws.Range("H139").Value = 69700
$ws.Range("J139").Value = 69700
$ws.Range("L139").Value = 69700
$ws.Range("N139").Value = -79980

# Row 140 (hunk 4)
$ws.Range("H140").Value = 77600
$ws.Range("J140").Value = 77600
$ws.Range("L140").Value = 77600
$ws.Range("N140").Value = -87960

# Row 141 (hunk 5)
$ws.Range("H141").Value = 850872.9399999999
$ws.Range("I141").Value = 1122071.9
$ws.Range("J141").Value = 3376.125
$ws.Range("K141").Value = 3366215.7
$ws.Range("L141").Value = 10128.375
$ws.Range("M141").Value = -3361035.7
$ws.Range("N141").Value = -20488.375

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (hunk 6)
$ws.Range("H32").Value = 5083.048
$ws.Range("I32").Value = 4638.1665
$ws.Range("J32").Value = 6195.25
$ws.Range("K32").Value = 4638.1665
$ws.Range("L32").Value = 6195.25
$ws.Range("M32").Value = -4351.1665
$ws.Range("N32").Value = -6769.25

# Row 61 (hunk 7)
$ws.Range("H61").Value = 2273.1667
$ws.Range("I61").Value = 871.625
$ws.Range("K61").Value = 871.625
$ws.Range("M61").Value = -659.625

# Row 74 (hunk 8)
$ws.Range("H74").Value = 1150.8077
$ws.Range("I74").Value = 784.1667
$ws.Range("K74").Value = 784.1667
$ws.Range("M74").Value = 89.83330000000001

# Row 77 (hunk 9)
$ws.Range("H77").Value = 1150.8077
$ws.Range("I77").Value = 784.1667
$ws.Range("K77").Value = 3920.8335
$ws.Range("M77").Value = 447.1665000000003

# Row 132 (hunk 10)
$ws.Range("H132").Value = 1390.9556
$ws.Range("I132").Value = 1125.5385
$ws.Range("J132").Value = 3116.1667
$ws.Range("K132").Value = 3376.6155
$ws.Range("L132").Value = 9348.500100000001
$ws.Range("M132").Value = -846.6155000000003
$ws.Range("N132").Value = -14408.5001

# Row 136 (hunk 11)
$ws.Range("H136").Value = 2273.1667
$ws.Range("I136").Value = 871.625
$ws.Range("K136").Value = 2614.875
$ws.Range("M136").Value = -64.875

$ws = $wb.Worksheets.Item("BSM")
# Row 99 (hunk 12)
$ws.Range("H99").Value = 1049.5
$ws.Range("I99").Value = 1100
$ws.Range("K99").Value = 1100
$ws.Range("M99").Value = 398

# Row 107 (hunk 13)
$ws.Range("H107").Value = 3389.5715
$ws.Range("I107").Value = 3458
$ws.Range("J107").Value = 2500
$ws.Range("K107").Value = 3458
$ws.Range("L107").Value = 2500
$ws.Range("M107").Value = -1538
$ws.Range("N107").Value = -6340

# Row 134 (hunk 14)
$ws.Range("H134").Value = 7493.9688
$ws.Range("I134").Value = 7872.6895
$ws.Range("K134").Value = 23618.0685
$ws.Range("M134").Value = -21083.0685

$ws = $wb.Worksheets.Item("CRP")
# Row 58 (hunk 15)
$ws.Range("H58").Value = 2558922.8
$ws.Range("I58").Value = 3953887.2
$ws.Range("K58").Value = 3953887.2
$ws.Range("M58").Value = -3953684.2

# Row 62 (hunk 16)
$ws.Range("H62").Value = 5684.3335
$ws.Range("I62").Value = 6279
$ws.Range("J62").Value = 3603
$ws.Range("K62").Value = 6279
$ws.Range("L62").Value = 3603
$ws.Range("M62").Value = -5655
$ws.Range("N62").Value = -4851

# Row 65 (hunk 17)
$ws.Range("H65").Value = 5684.3335
$ws.Range("I65").Value = 6279
$ws.Range("J65").Value = 3603
$ws.Range("K65").Value = 31395
$ws.Range("L65").Value = 18015
$ws.Range("M65").Value = -28275
$ws.Range("N65").Value = -24255

# Row 86 (hunk 18)
$ws.Range("H86").Value = 90911496
$ws.Range("I86").Value = 111113470
$ws.Range("K86").Value = 111113470
$ws.Range("M86").Value = -111112347

# Row 89 (hunk 19)
$ws.Range("H89").Value = 90911496
$ws.Range("I89").Value = 111113470
$ws.Range("K89").Value = 555567350
$ws.Range("M89").Value = -555561734

# Row 99 (hunk 20)
$ws.Range("H99").Value = 3582.818
$ws.Range("I99").Value = 2799.625
$ws.Range("K99").Value = 2799.625
$ws.Range("M99").Value = -1301.625

# Row 122 (hunk 21)
$ws.Range("H122").Value = 1446.9474
$ws.Range("I122").Value = 1036
$ws.Range("J122").Value = 2337.3333
$ws.Range("K122").Value = 3108
$ws.Range("L122").Value = 7011.999899999999
$ws.Range("M122").Value = -658
$ws.Range("N122").Value = -11911.9999

# Row 126 (hunk 22)
$ws.Range("H126").Value = 3582.818
$ws.Range("I126").Value = 2799.625
$ws.Range("K126").Value = 8398.875
$ws.Range("M126").Value = -5928.875

# Row 132 (hunk 23)
$ws.Range("H132").Value = 2307.3684
$ws.Range("I132").Value = 1675.1666
$ws.Range("J132").Value = 4678.125
$ws.Range("K132").Value = 5025.4998
$ws.Range("L132").Value = 14034.375
$ws.Range("M132").Value = -2495.4998
$ws.Range("N132").Value = -19094.375

# Row 133 (hunk 24)
$ws.Range("H133").Value = 60000
$ws.Range("J133").Value = 60000
$ws.Range("L133").Value = 60000
$ws.Range("N133").Value = -65060

# Row 134 (hunk 25)
$ws.Range("H134").Value = 1505.84
$ws.Range("I134").Value = 838.1795
$ws.Range("J134").Value = 3873
$ws.Range("K134").Value = 2514.5385
$ws.Range("L134").Value = 11619
$ws.Range("M134").Value = 20.46150000000034
$ws.Range("N134").Value = -16689

# Row 135 (hunk 26)
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# Row 136 (hunk 27)
$ws.Range("H136").Value = 2558922.8
$ws.Range("I136").Value = 3953887.2
$ws.Range("K136").Value = 11861661.6
$ws.Range("M136").Value = -11859111.6

$ws = $wb.Worksheets.Item("CUL")
# Row 33 (hunk 28)
$ws.Range("H33").Value = 221.8
$ws.Range("J33").Value = 380.5
$ws.Range("L33").Value = 2283
$ws.Range("N33").Value = -2849

# Row 100 (hunk 29)
$ws.Range("H100").Value = 3105

# Row 131 (hunk 30)
$ws.Range("H131").Value = 8633783
$ws.Range("J131").Value = 13783.728
$ws.Range("L131").Value = 41351.18399999999
$ws.Range("N131").Value = -51431.18399999999

$ws = $wb.Worksheets.Item("GSM")
# Row 11 (hunk 31)
$ws.Range("H11").Value = 5472718.5
$ws.Range("I11").Value = 2336480.5
$ws.Range("J11").Value = 9589031
$ws.Range("K11").Value = 2336480.5
$ws.Range("L11").Value = 9589031
$ws.Range("M11").Value = -2336341.5
$ws.Range("N11").Value = -9589309

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (hunk 32)
$ws.Range("H7").Value = 4680.1816
$ws.Range("I7").Value = 2685.5
$ws.Range("J7").Value = 9999.333000000001
$ws.Range("K7").Value = 2685.5
$ws.Range("L7").Value = 9999.333000000001
$ws.Range("M7").Value = -2573.5
$ws.Range("N7").Value = -10223.333

# Row 61 (hunk 33)
$ws.Range("H61").Value = 2517.4119
$ws.Range("I61").Value = 1985.5
$ws.Range("K61").Value = 1985.5
$ws.Range("M61").Value = -1783.5

# Row 113 (hunk 34)
$ws.Range("H113").Value = 2517.4119
$ws.Range("I113").Value = 1985.5
$ws.Range("K113").Value = 1985.5
$ws.Range("M113").Value = 184.5

# Row 126 (hunk 35)
$ws.Range("H126").Value = 4680.1816
$ws.Range("I126").Value = 2685.5
$ws.Range("J126").Value = 9999.333000000001
$ws.Range("K126").Value = 8056.5
$ws.Range("L126").Value = 29997.999
$ws.Range("M126").Value = -5586.5
$ws.Range("N126").Value = -34937.999

$ws = $wb.Worksheets.Item("WVR")
# Row 28 (hunk 36)
$ws.Range("H28").Value = 9612.5
$ws.Range("I28").Value = 8500
$ws.Range("J28").Value = 9983.333000000001
$ws.Range("K28").Value = 8500
$ws.Range("L28").Value = 9983.333000000001
$ws.Range("M28").Value = -8152
$ws.Range("N28").Value = -10679.333

# Row 124 (hunk 37)
$ws.Range("H124").Value = 19800
$ws.Range("J124").Value = 19800
$ws.Range("L124").Value = 19800
$ws.Range("N124").Value = -29620

# Row 132 (hunk 38)
$ws.Range("H132").Value = 808.48
$ws.Range("I132").Value = 338.10526
$ws.Range("J132").Value = 2298
$ws.Range("K132").Value = 1014.31578
$ws.Range("L132").Value = 6894
$ws.Range("M132").Value = 1515.68422
$ws.Range("N132").Value = -11954
